$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2962
$ws.Range("J17").Value = 2962
$ws.Range("L17").Value = 8886
$ws.Range("N17").Value = -9222
$ws.Range("H32").Value = 2525.6667
$ws.Range("I32").Value = 2047
$ws.Range("J32").Value = 2867.5715
$ws.Range("K32").Value = 2047
$ws.Range("L32").Value = 2867.5715
$ws.Range("M32").Value = -1721
$ws.Range("N32").Value = -3519.5715
$ws.Range("H74").Value = 5970
$ws.Range("I74").Value = 5411
$ws.Range("K74").Value = 5411
$ws.Range("M74").Value = -4475
$ws.Range("H77").Value = 5970
$ws.Range("I77").Value = 5411
$ws.Range("K77").Value = 27055
$ws.Range("M77").Value = -22375
$ws.Range("H80").Value = 1422.2632
$ws.Range("I80").Value = 1034.1111
$ws.Range("K80").Value = 3102.3333
$ws.Range("M80").Value = -2104.3333
$ws.Range("H83").Value = 1422.2632
$ws.Range("I83").Value = 1034.1111
$ws.Range("K83").Value = 9306.999900000001
$ws.Range("M83").Value = -4314.999900000001
$ws.Range("H100").Value = 2354.5557
$ws.Range("I100").Value = 2354.5557
$ws.Range("K100").Value = 2354.5557
$ws.Range("M100").Value = -1813.5557
$ws.Range("H103").Value = 605.5
$ws.Range("I103").Value = 679.3333
$ws.Range("J103").Value = 494.75
$ws.Range("K103").Value = 2037.9999
$ws.Range("L103").Value = 1484.25
$ws.Range("M103").Value = -1451.9999
$ws.Range("N103").Value = -2656.25
$ws.Range("H138").Value = 1229.6666
$ws.Range("I138").Value = 1202.5294
$ws.Range("K138").Value = 3607.5882
$ws.Range("M138").Value = 1532.4118

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 2003
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("H110").Value = 6998.25
$ws.Range("I110").Value = 8206.625
$ws.Range("J110").Value = 5387.0835
$ws.Range("K110").Value = 8206.625
$ws.Range("L110").Value = 5387.0835
$ws.Range("M110").Value = -6161.625
$ws.Range("N110").Value = -9477.083500000001
$ws.Range("H132").Value = 3967.9524
$ws.Range("I132").Value = 3293.9033
$ws.Range("J132").Value = 5867.5454
$ws.Range("K132").Value = 9881.7099
$ws.Range("L132").Value = 17602.6362
$ws.Range("M132").Value = -7351.7099
$ws.Range("N132").Value = -22662.6362
$ws.Range("N12").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 4015.875
$ws.Range("I64").Value = 2959
$ws.Range("K64").Value = 2959
$ws.Range("M64").Value = -2734
$ws.Range("H67").Value = 4015.875
$ws.Range("I67").Value = 2959
$ws.Range("K67").Value = 2959
$ws.Range("M67").Value = -2179
$ws.Range("H86").Value = 11217.556
$ws.Range("I86").Value = 7036.2
$ws.Range("K86").Value = 7036.2
$ws.Range("M86").Value = -5913.2
$ws.Range("H89").Value = 11217.556
$ws.Range("I89").Value = 7036.2
$ws.Range("K89").Value = 35181
$ws.Range("M89").Value = -29565
$ws.Range("H105").Value = 2549.8
$ws.Range("I105").Value = 2437.25
$ws.Range("K105").Value = 2437.25
$ws.Range("M105").Value = -690.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 2851.875
$ws.Range("I10").Value = 1007.2727
$ws.Range("J10").Value = 6910
$ws.Range("K10").Value = 1007.2727
$ws.Range("L10").Value = 6910
$ws.Range("M10").Value = -868.2727
$ws.Range("N10").Value = -7188
$ws.Range("H12").Value = 5116.6665
$ws.Range("I12").Value = 2675
$ws.Range("J12").Value = 10000
$ws.Range("K12").Value = 2675
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = -2505
$ws.Range("N12").Value = -10340
$ws.Range("H16").Value = 1779.4
$ws.Range("I16").Value = 1477.5
$ws.Range("J16").Value = 2987
$ws.Range("K16").Value = 1477.5
$ws.Range("L16").Value = 2987
$ws.Range("M16").Value = -1190.5
$ws.Range("N16").Value = -3561
$ws.Range("H59").Value = 69999
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("H97").Value = 70000
$ws.Range("I97").Value = 70000
$ws.Range("K97").Value = 70000
$ws.Range("M97").Value = -69009
$ws.Range("H107").Value = 660.3
$ws.Range("I107").Value = 544.8333
$ws.Range("K107").Value = 544.8333
$ws.Range("M107").Value = 1375.1667
$ws.Range("H113").Value = 1779.4
$ws.Range("I113").Value = 1477.5
$ws.Range("J113").Value = 2987
$ws.Range("K113").Value = 1477.5
$ws.Range("L113").Value = 2987
$ws.Range("M113").Value = 692.5
$ws.Range("N113").Value = -7327
$ws.Range("M59").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 2999.6667
$ws.Range("J49").Value = 2250
$ws.Range("L49").Value = 6750
$ws.Range("N49").Value = -7062
$ws.Range("H92").Value = 402.0909
$ws.Range("I92").Value = 403.77777
$ws.Range("J92").Value = 394.5
$ws.Range("K92").Value = 1211.33331
$ws.Range("L92").Value = 1183.5
$ws.Range("M92").Value = 36.66669000000002
$ws.Range("N92").Value = -3679.5
$ws.Range("H130").Value = 1999.6666

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1811.2142
$ws.Range("I97").Value = 1896.125
$ws.Range("J97").Value = 1698
$ws.Range("K97").Value = 1896.125
$ws.Range("L97").Value = 1698
$ws.Range("M97").Value = -1400.125
$ws.Range("N97").Value = -2690
$ws.Range("H100").Value = 25000
$ws.Range("J100").Value = 25000
$ws.Range("L100").Value = 25000
$ws.Range("N100").Value = -27164
$ws.Range("H122").Value = 2239.8
$ws.Range("I122").Value = 2221.2144
$ws.Range("K122").Value = 6663.6432
$ws.Range("M122").Value = -4213.6432
$ws.Range("H132").Value = 2205.7932
$ws.Range("I132").Value = 2159.5833
$ws.Range("J132").Value = 2427.6
$ws.Range("K132").Value = 6478.749899999999
$ws.Range("L132").Value = 7282.799999999999
$ws.Range("M132").Value = -3948.749899999999
$ws.Range("N132").Value = -12342.8
$ws.Range("H134").Value = 43899
$ws.Range("J134").Value = 43899
$ws.Range("L134").Value = 131697
$ws.Range("N134").Value = -136767

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 15003.5
$ws.Range("J11").Value = 15003.5
$ws.Range("L11").Value = 15003.5
$ws.Range("N11").Value = -15283.5
$ws.Range("H16").Value = 375
$ws.Range("I16").Value = 375
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 375
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = -205
$ws.Range("H18").Value = 4999
$ws.Range("I18").Value = 4999
$ws.Range("K18").Value = 4999
$ws.Range("M18").Value = -4827
$ws.Range("H100").Value = 1997.8
$ws.Range("I100").Value = 1997.5
$ws.Range("K100").Value = 1997.5
$ws.Range("M100").Value = -1456.5
$ws.Range("H101").Value = 6965.6
$ws.Range("J101").Value = 6965.6
$ws.Range("L101").Value = 6965.6
$ws.Range("N101").Value = -13455.6
$ws.Range("M16").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 6402.8335
$ws.Range("I96").Value = 6192
$ws.Range("K96").Value = 6192
$ws.Range("M96").Value = -4819
$ws.Range("H100").Value = 258
$ws.Range("I100").Value = 258
$ws.Range("K100").Value = 516
$ws.Range("M100").Value = 25
$ws.Range("H113").Value = 407.54544
$ws.Range("I113").Value = 294.7143
$ws.Range("J113").Value = 605
$ws.Range("K113").Value = 884.1428999999999
$ws.Range("L113").Value = 1815
$ws.Range("M113").Value = 1285.8571
$ws.Range("N113").Value = -6155
$ws.Range("H123").Value = 60000
$ws.Range("J123").Value = 60000
$ws.Range("L123").Value = 60000
$ws.Range("N123").Value = -69800
$ws.Range("H133").Value = 89999.664
$ws.Range("J133").Value = 89999.664
$ws.Range("L133").Value = 89999.664
$ws.Range("N133").Value = -100119.664
